$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old B3 value (removes "EB000016" from the shared-string table)
$ws.Range("B3").Value = ""

# Update B2 in place (EB000013 -> EB001250)
$ws.Range("B2").Value = "EB001250"

# New cells, written in column-within-row order to build the shared-string
# table in the same append order as the target file.
$ws.Range("D1").Value = "productNum"
$ws.Range("F1").Value = "productDesc"
$ws.Range("E1").Value = "productOpt"
$ws.Range("D2").Value = "H1K92A3"
$ws.Range("F2").Value = "HPE 3Y Proactive Care 24x7 SVC"
$ws.Range("F3").Value = "HP Install DL36x(p) Service"
$ws.Range("D3").Value = "U4506E"

$ws.Range("D7").Select()
